# PRJ0018286-CF Industry Group Changes TECH
# - Rename "JobType" sheet to "FilterRecord"
# - Update header text on FilterRecord sheet and add new "TECH - Technology" row
# - Add "Engagements" to ReportOption sheet
# - Add "Industry Group" to Filter sheet
# - Filter sheet becomes the active/selected sheet

$wb = $excel.ActiveWorkbook

$wsFilterRecord = $wb.Worksheets.Item(1)
$wsReportOption = $wb.Worksheets.Item(2)
$wsFilter       = $wb.Worksheets.Item(3)

# --- Sheet1: JobType -> FilterRecord ---
$wsFilterRecord.Name = "FilterRecord"
$wsFilterRecord.Range("A1").Value = "FilterRecordValue"
$wsFilterRecord.Range("A3").Value = "TECH - Technology"
$wsFilterRecord.Columns.Item(1).ColumnWidth = 33.7

# --- Sheet2: ReportOption gains "Engagements" ---
$wsReportOption.Range("A3").Value = "Engagements"

# --- Sheet3: Filter gains "Industry Group" ---
$wsFilter.Range("A3").Value = "Industry Group"
$wsFilter.Columns.Item(1).ColumnWidth = 12.2

# --- Selections on each sheet ---
$wsFilterRecord.Range("C21").Select() | Out-Null
$wsReportOption.Range("B8").Select() | Out-Null
$wsFilter.Range("E24").Select() | Out-Null

# --- Filter sheet is the active/selected tab ---
$wsFilter.Activate()
